# Refresh Marketboard-derived Leve profit figures (H:N) across sheets.
# Values below come from an external scheduled price-sync run; only
# numeric cells in columns H-N on the listed rows are touched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19 (Leve Item ID 7015)
$ws.Range("H19").Value = 174.83333
$ws.Range("I19").Value = 174.83333
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 174.83333
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 0.1666700000000105
$ws.Range("N19").ClearContents()

# Row 40 (Leve Item ID 5505)
$ws.Range("H40").Value = 3566.2222
$ws.Range("J40").Value = 3637.125
$ws.Range("L40").Value = 3637.125
$ws.Range("N40").Value = -3987.125

# Row 141 (Leve Item ID 44161)
$ws.Range("H141").Value = 4618.5
$ws.Range("I141").Value = 4618.5
$ws.Range("K141").Value = 13855.5
$ws.Range("M141").Value = -8675.5

$ws = $wb.Worksheets.Item("ARM")
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 2733.7778
$ws.Range("J45").Value = 3497.5
$ws.Range("L45").Value = 3497.5
$ws.Range("N45").Value = -4251.5

# Row 54 (Leve Item ID 2817)
$ws.Range("H54").Value = 10499.5
$ws.Range("I54").Value = 1000
$ws.Range("J54").Value = 19999
$ws.Range("K54").Value = 1000
$ws.Range("L54").Value = 19999
$ws.Range("M54").Value = -231
$ws.Range("N54").Value = -21537

# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 2081.75
$ws.Range("I61").Value = 2081.75
$ws.Range("K61").Value = 2081.75
$ws.Range("M61").Value = -1869.75

# Row 75 (Leve Item ID 10714)
$ws.Range("H75").Value = 142499
$ws.Range("J75").Value = 142499
$ws.Range("L75").Value = 142499
$ws.Range("N75").Value = -144247

# Row 78 (Leve Item ID 10714)
$ws.Range("H78").Value = 142499
$ws.Range("J78").Value = 142499
$ws.Range("L78").Value = 427497
$ws.Range("N78").Value = -436233

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 3793.8518
$ws.Range("I132").Value = 3516.6924
$ws.Range("K132").Value = 10550.0772
$ws.Range("M132").Value = -8020.0772

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 2081.75
$ws.Range("I136").Value = 2081.75
$ws.Range("K136").Value = 6245.25
$ws.Range("M136").Value = -3695.25

$ws = $wb.Worksheets.Item("BSM")
# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 2178.8333
$ws.Range("I86").Value = 2336.3845
$ws.Range("J86").Value = 1769.2
$ws.Range("K86").Value = 2336.3845
$ws.Range("L86").Value = 1769.2
$ws.Range("M86").Value = -1213.3845
$ws.Range("N86").Value = -4015.2

# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 2178.8333
$ws.Range("I89").Value = 2336.3845
$ws.Range("J89").Value = 1769.2
$ws.Range("K89").Value = 11681.9225
$ws.Range("L89").Value = 8846
$ws.Range("M89").Value = -6065.922500000001
$ws.Range("N89").Value = -20078

$ws = $wb.Worksheets.Item("CRP")
# Row 29 (Leve Item ID 2408)
$ws.Range("H29").Value = 1033
$ws.Range("J29").Value = 1033
$ws.Range("L29").Value = 1033
$ws.Range("N29").Value = -1619

# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 3513.25
$ws.Range("J31").Value = 3184.5
$ws.Range("L31").Value = 3184.5
$ws.Range("N31").Value = -3774.5

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 3513.25
$ws.Range("J34").Value = 3184.5
$ws.Range("L34").Value = 3184.5
$ws.Range("N34").Value = -3588.5

# Row 63 (Leve Item ID 10604)
$ws.Range("H63").Value = 69999.5
$ws.Range("J63").Value = 69999.5
$ws.Range("L63").Value = 69999.5
$ws.Range("N63").Value = -71371.5

# Row 66 (Leve Item ID 10604)
$ws.Range("H66").Value = 69999.5
$ws.Range("J66").Value = 69999.5
$ws.Range("L66").Value = 209998.5
$ws.Range("N66").Value = -216862.5

# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 3425
$ws.Range("I99").Value = 1900
$ws.Range("K99").Value = 1900
$ws.Range("M99").Value = -402

# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 3425
$ws.Range("I126").Value = 1900
$ws.Range("K126").Value = 5700
$ws.Range("M126").Value = -3230

# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 2710.6667
$ws.Range("I134").Value = 2589.5386
$ws.Range("J134").Value = 3498
$ws.Range("K134").Value = 7768.6158
$ws.Range("L134").Value = 10494
$ws.Range("M134").Value = -5233.6158
$ws.Range("N134").Value = -15564

$ws = $wb.Worksheets.Item("CUL")
# Row 7 (Leve Item ID 4728)
$ws.Range("H7").Value = 480
$ws.Range("J7").Value = 400
$ws.Range("L7").Value = 1200
$ws.Range("N7").Value = -1424

# Row 23 (Leve Item ID 4858)
$ws.Range("H23").Value = 799.6667
$ws.Range("I23").Value = 799
$ws.Range("J23").Value = 799.75
$ws.Range("K23").Value = 2397
$ws.Range("L23").Value = 2399.25
$ws.Range("M23").Value = -2162
$ws.Range("N23").Value = -2869.25

# Row 55 (Leve Item ID 4733)
$ws.Range("H55").Value = 8627.700000000001
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 8627.700000000001
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 25883.1
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -26237.1

# Row 69 (Leve Item ID 12850)
$ws.Range("H69").Value = 2360
$ws.Range("I69").Value = 2500
$ws.Range("J69").Value = 2325
$ws.Range("K69").Value = 7500
$ws.Range("L69").Value = 6975
$ws.Range("M69").Value = -6689
$ws.Range("N69").Value = -8597

# Row 72 (Leve Item ID 12850)
$ws.Range("H72").Value = 2360
$ws.Range("I72").Value = 2500
$ws.Range("J72").Value = 2325
$ws.Range("K72").Value = 22500
$ws.Range("L72").Value = 20925
$ws.Range("M72").Value = -18444
$ws.Range("N72").Value = -29037

$ws = $wb.Worksheets.Item("GSM")
# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 60710.7
$ws.Range("I122").Value = 62138.375
$ws.Range("K122").Value = 186415.125
$ws.Range("M122").Value = -183965.125

# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 1949.5
$ws.Range("I126").Value = 1949.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5848.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3378.5
$ws.Range("N126").ClearContents()

# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 1295.8
$ws.Range("I132").Value = 994.75
$ws.Range("K132").Value = 2984.25
$ws.Range("M132").Value = -454.25

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 2692.2222
$ws.Range("I16").Value = 2841.25
$ws.Range("J16").Value = 1500
$ws.Range("K16").Value = 2841.25
$ws.Range("L16").Value = 1500
$ws.Range("M16").Value = -2671.25
$ws.Range("N16").Value = -1840

# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 889.2727
$ws.Range("I55").Value = 910.25
$ws.Range("K55").Value = 910.25
$ws.Range("M55").Value = -737.25

# Row 68 (Leve Item ID 12563)
$ws.Range("H68").Value = 3297.3333
$ws.Range("I68").Value = 3056.8
$ws.Range("K68").Value = 3056.8
$ws.Range("M68").Value = -2307.8

# Row 71 (Leve Item ID 12563)
$ws.Range("H71").Value = 3297.3333
$ws.Range("I71").Value = 3056.8
$ws.Range("K71").Value = 15284
$ws.Range("M71").Value = -11540

# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 3208.8
$ws.Range("I122").Value = 3065.3333
$ws.Range("K122").Value = 9195.999899999999
$ws.Range("M122").Value = -6745.999899999999

$ws = $wb.Worksheets.Item("WVR")
# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 2574.5
$ws.Range("I126").Value = 2599.3333
$ws.Range("K126").Value = 7797.999899999999
$ws.Range("M126").Value = -5327.999899999999
